# Apply updated "dSF" (column F) values for specific rows in the save-data
# sheet. These correspond to a re-pull of the underlying data, so only the
# final/"F" snapshot column values change while the initial ("E"/dS0) values
# are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 4
    9  = 3
    16 = 3
    19 = 0
    24 = -3
    26 = 2
    32 = 1
    38 = 2
    41 = -3
    49 = -6
    51 = -1
    52 = -6
    53 = -1
    54 = -1
    56 = 2
    58 = 0
    61 = -3
    65 = 1
    66 = -2
    67 = -2
    69 = -1
    72 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
